$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 60
$ws.Range("E2").Value = "'43,403,887.00"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'10.55"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 25.071
$ws.Range("D3").Value = 58
$ws.Range("E3").Value = "'43,873,397.00"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'10.81"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 27.951
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = "'28,345,210.00"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'6.70"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 21.809
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = "'4.68"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = "'5.43"
$ws.Range("L4").Style = "Normal"
$ws.Range("D5").Value = 65
$ws.Range("E5").Value = "'35,130,384.00"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'9.07"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 24.912
$ws.Range("J5").Value = 19
$ws.Range("K5").Value = "'4.86"
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").Value = "'5.69"
$ws.Range("L5").Style = "Normal"
$ws.Range("H6").Value = 19.429
$ws.Range("D7").Value = 76
$ws.Range("E7").Value = "'45,216,429.00"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'10.40"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = 14.664
$ws.Range("J7").Value = 18
$ws.Range("K7").Value = "'4.32"
$ws.Range("K7").Style = "Normal"
$ws.Range("L7").Value = "'4.90"
$ws.Range("L7").Style = "Normal"
$ws.Range("D8").Value = 63
$ws.Range("E8").Value = "'48,832,401.00"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'12.08"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = 30.324
$ws.Range("K8").Value = "'3.61"
$ws.Range("K8").Style = "Normal"
$ws.Range("D9").Value = 49
$ws.Range("E9").Value = "'41,504,223.00"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'9.18"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = 16.432
$ws.Range("D10").Value = 38
$ws.Range("E10").Value = "'24,747,676.00"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'6.83"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = 17.897
$ws.Range("H11").Value = 16.228
$ws.Range("E12").Value = "'34,024,953.00"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'8.39"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = 25.081
$ws.Range("J12").Value = 19
$ws.Range("K12").Value = "'4.96"
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = "'5.15"
$ws.Range("L12").Style = "Normal"
$ws.Range("H13").Value = 20.81
$ws.Range("D14").Value = 65
$ws.Range("E14").Value = "'56,632,594.00"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'14.03"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = 49
$ws.Range("H15").Value = 16.532
$ws.Range("D16").Value = 56
$ws.Range("E16").Value = "'37,171,494.00"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'8.99"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = 13.184
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = "'2.78"
$ws.Range("K16").Style = "Normal"
$ws.Range("L16").Value = "'3.80"
$ws.Range("L16").Style = "Normal"
$ws.Range("D18").Value = 38
$ws.Range("E18").Value = "'29,826,602.00"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'8.51"
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value = 16.669
$ws.Range("K18").Value = "'3.67"
$ws.Range("K18").Style = "Normal"

$ws.Name = "repayment_20250901_20250923"
